# Update schedule data for Línea 141 (scrape refresh at 03:40:31)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 03:40:31"
$ws1.Range("A3").Value = "Total filas: 7"

$rows1 = @(
    @("03:40:31", "03:48", "14_ABASTO",      8,   "LP1912"),
    @("03:40:31", "04:01", "81_EL PELIGRO",  21,  "LP1912"),
    @("03:40:31", "04:46", "215A_EL PATO",   66,  "LP1912"),
    @("03:40:31", "04:53", "11_ETCHEVERRY",  73,  "LP1912"),
    @("03:40:31", "05:16", "17_ROMERO",      96,  "LP1912"),
    @("03:40:31", "05:22", "23_HERNANDEZ",   102, "LP1912"),
    @("03:40:31", "05:34", "215B_EL PATO",   114, "LP1912")
)

$r = 6
foreach ($row in $rows1) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 03:40:31"

$rows2 = @(
    @("03:40:31", "04:46", "215A_EL PATO", 66,  "LP1912"),
    @("03:40:31", "05:34", "215B_EL PATO", 114, "LP1912")
)

$r = 6
foreach ($row in $rows2) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 03:40:31"

Write-Output "Update complete"
